$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped entirely from the dataset
# (row 28 = "SC 92", row 26 = "RM 232"). Delete the higher index first
# so the lower index still refers to the correct row afterwards.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two row deletions, apply the individual cell edits
# (values newly imputed/filled in, or newly blanked out) using the
# resulting (post-delete) row numbers.

$ws.Range("C2").Value = 14.9
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = -6.4
$ws.Range("E5").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("C12").Value = 12.5
$ws.Range("C14").Value = ""
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("E23").Value = -7

$ws.Range("E27").Value = ""
$ws.Range("E29").Value = -6.8
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("B32").Value = ""
$ws.Range("C33").Value = 10.4
